$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5882352941176471
$ws.Range("D2").Value = 0.6741573033707866
$ws.Range("B3").Value = 0.8095238095238095
$ws.Range("C3").Value = 0.6181818181818182
$ws.Range("D3").Value = 0.7010309278350516
$ws.Range("B4").Value = 0.6881720430107527
$ws.Range("C4").Value = 0.6881720430107527
$ws.Range("D4").Value = 0.6881720430107527
$ws.Range("E4").Value = 0.6881720430107527
$ws.Range("B5").Value = 0.6988795518207283
$ws.Range("C5").Value = 0.7038277511961722
$ws.Range("D5").Value = 0.6875941156029191
$ws.Range("B6").Value = 0.719104846239571
$ws.Range("C6").Value = 0.6881720430107527
$ws.Range("D6").Value = 0.6900503070862122
$ws.Range("B7").Value = 0.6279069767441861
$ws.Range("C7").Value = 0.7105263157894737
$ws.Range("D7").Value = 0.6666666666666666
$ws.Range("B8").Value = 0.78
$ws.Range("C8").Value = 0.7090909090909091
$ws.Range("D8").Value = 0.7428571428571428
$ws.Range("B9").Value = 0.7096774193548387
$ws.Range("C9").Value = 0.7096774193548387
$ws.Range("D9").Value = 0.7096774193548387
$ws.Range("E9").Value = 0.7096774193548387
$ws.Range("B10").Value = 0.703953488372093
$ws.Range("C10").Value = 0.7098086124401914
$ws.Range("D10").Value = 0.7047619047619047
$ws.Range("B11").Value = 0.7178544636159039
$ws.Range("C11").Value = 0.7096774193548387
$ws.Range("D11").Value = 0.7117255504352279
$ws.Range("B12").Value = 0.7941176470588235
$ws.Range("C12").Value = 0.7105263157894737
$ws.Range("D12").Value = 0.7499999999999999
$ws.Range("B13").Value = 0.8135593220338984
$ws.Range("C13").Value = 0.8727272727272727
$ws.Range("D13").Value = 0.8421052631578948
$ws.Range("B14").Value = 0.8064516129032258
$ws.Range("C14").Value = 0.8064516129032258
$ws.Range("D14").Value = 0.8064516129032258
$ws.Range("E14").Value = 0.8064516129032258
$ws.Range("B15").Value = 0.8038384845463609
$ws.Range("C15").Value = 0.7916267942583732
$ws.Range("D15").Value = 0.7960526315789473
$ws.Range("B16").Value = 0.8056154118290291
$ws.Range("C16").Value = 0.8064516129032258
$ws.Range("D16").Value = 0.8044708545557442
$ws.Range("B17").Value = 0.6666666666666666
$ws.Range("C17").Value = 0.7368421052631579
$ws.Range("D17").Value = 0.7
$ws.Range("B18").Value = 0.803921568627451
$ws.Range("C18").Value = 0.7454545454545455
$ws.Range("D18").Value = 0.7735849056603775
$ws.Range("B19").Value = 0.7419354838709677
$ws.Range("C19").Value = 0.7419354838709677
$ws.Range("D19").Value = 0.7419354838709677
$ws.Range("E19").Value = 0.7419354838709677
$ws.Range("B20").Value = 0.7352941176470589
$ws.Range("C20").Value = 0.7411483253588517
$ws.Range("D20").Value = 0.7367924528301888
$ws.Range("B21").Value = 0.7478389205144423
$ws.Range("C21").Value = 0.7419354838709677
$ws.Range("D21").Value = 0.7435179549604384
$ws.Range("B22").Value = 0.6923076923076923
$ws.Range("C22").Value = 0.7105263157894737
$ws.Range("D22").Value = 0.7012987012987013
$ws.Range("B23").Value = 0.7962962962962963
$ws.Range("C23").Value = 0.7818181818181819
$ws.Range("D23").Value = 0.7889908256880734
$ws.Range("B24").Value = 0.7526881720430108
$ws.Range("C24").Value = 0.7526881720430108
$ws.Range("D24").Value = 0.7526881720430108
$ws.Range("E24").Value = 0.7526881720430108
$ws.Range("B25").Value = 0.7443019943019943
$ws.Range("C25").Value = 0.7461722488038278
$ws.Range("D25").Value = 0.7451447634933874
$ws.Range("B26").Value = 0.7538063290751462
$ws.Range("C26").Value = 0.7526881720430108
$ws.Range("D26").Value = 0.7531596350773623